$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data rows (row 2..9), columns A-G
# (numbers are written in plain decimal form since the interpreter does not
#  accept scientific-notation literals like 1e-12)
$data = @(
    @("even_MAG-GUT1861.fa",  0.03898443179991121,    0.9609849467234637,   0.00003062147662507972,    0.9609849467234637,   "o__Enterobacterales", "o__Enterobacterales"),
    @("even_MAG-GUT3427.fa",  0.04234714895744301,    0.9576073625115518,   0.00004548853100520531,    0.9576073625115518,   "o__Enterobacterales", "o__Enterobacterales"),
    @("even_MAG-GUT3435.fa",  0.04044130900443254,    0.9594889110945014,   0.0000697799010660079,     0.9594889110945014,   "o__Enterobacterales", "o__Enterobacterales"),
    @("even_MAG-GUT43440.fa", 0.000000000005477520756733674, 0.9945728552818769, 0.005427144712645665, 0.9945728552818769,   "o__Enterobacterales", "o__Enterobacterales"),
    @("even_MAG-GUT8151.fa",  0.9968757625693284,     0.002802703324043471, 0.0003215341066280073,     0.9968757625693284,   "o__Aeromonadales",    "o__Aeromonadales"),
    @("even_MAG-GUT91566.fa", 0.05172531570058501,    0.9480379032992007,   0.000236781000214267,      0.9480379032992007,   "o__Enterobacterales", "o__Enterobacterales"),
    @("even_MAG-GUT91702.fa", 0.000000000005060160787590307, 0.9878933824326019, 0.0121066175623379,   0.9878933824326019,   "o__Enterobacterales", "o__Enterobacterales"),
    @("even_MAG-GUT91898.fa", 0.000000000001836383938701958, 0.9933893857680852, 0.006610614230078402, 0.9933893857680852,   "o__Enterobacterales", "o__Enterobacterales")
)

# Write all the values first (rows 2-9, columns A-G)
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}

# Apply the header-style formatting (bold, centered, thin border) to column A of every
# data row, matching the style used by the original rows (cell style index 1). Using
# copy/paste-special of formats (rather than assigning individual Font/Alignment
# properties) keeps every row mapped onto the same, single shared style.
$ws.Cells.Item(1, 1).Copy() | Out-Null
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false
